# "some minor alignment improved"
#
# Slide 3 has a group ("Group 103", id 104) holding two shapes:
#   - id 107 "Snip Diagonal Corner Rectangle 106"
#   - id 108 "TextBox 107" (three separate runs forming one quote)
#
# The edit ungroups them (flattening their transforms into the slide's
# coordinate space) and then nudges their position/size slightly, and
# collapses the textbox's three runs into a single run.

$EMU = 12700.0

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Locate the group shape by name (robust against index drift).
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Group 103") {
        $grp = $candidate
        break
    }
}

$items = $grp.Ungroup()

for ($i = 1; $i -le $items.Count; $i++) {
    $sh = $items.Item($i)

    if ($sh.Name -eq "Snip Diagonal Corner Rectangle 106") {
        # tiny epsilon nudges the float32 Left/Top past the EMU rounding edge
        $sh.Left = 7426306 / $EMU + 0.00001
        $sh.Top = 5260903 / $EMU + 0.00001
        $sh.Width = 1166672 / $EMU
        $sh.Height = 1487240 / $EMU
    }
    elseif ($sh.Name -eq "TextBox 107") {
        $sh.Left = 7402700 / $EMU
        $sh.Top = 5594196 / $EMU
        $sh.Width = 1224136 / $EMU
        $sh.Height = 769441 / $EMU
        $sh.TextFrame.TextRange.Text = "“Insomniac presents Volume Sundays featuring Zeds Dead…”"
    }
}
